$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value2 = 547.7778
$ws.Range("I33").Value2 = 528.75
$ws.Range("J33").Value2 = 700
$ws.Range("K33").Value2 = 528.75
$ws.Range("L33").Value2 = 700
$ws.Range("M33").Value2 = -299.75
$ws.Range("N33").Value2 = -1158

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value2 = 40335.5
$ws.Range("J105").Value2 = 40335.5
$ws.Range("L105").Value2 = 40335.5
$ws.Range("N105").Value2 = -47323.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value2 = 2104.4443
$ws.Range("J112").Value2 = 2190.5881
$ws.Range("L112").Value2 = 6571.7643
$ws.Range("N112").Value2 = -8787.764299999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value2 = 1568.8462
$ws.Range("I113").Value2 = 1000
$ws.Range("J113").Value2 = 1672.2727
$ws.Range("K113").Value2 = 1000
$ws.Range("L113").Value2 = 1672.2727
$ws.Range("M113").Value2 = 2254
$ws.Range("N113").Value2 = -8180.2727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 24391796
$ws.Range("I137").Value2 = 1063.3438
$ws.Range("J137").Value2 = 111114400
$ws.Range("K137").Value2 = 3190.0314
$ws.Range("L137").Value2 = 333343200
$ws.Range("M137").Value2 = -640.0314000000003
$ws.Range("N137").Value2 = -333348300

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value2 = 72000
$ws.Range("J140").Value2 = 72000
$ws.Range("L140").Value2 = 72000
$ws.Range("N140").Value2 = -82360

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value2 = 1862.7
$ws.Range("I141").Value2 = 1054.909
$ws.Range("K141").Value2 = 3164.727
$ws.Range("M141").Value2 = 2015.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 1062.6818
$ws.Range("I61").Value2 = 912.4375
$ws.Range("J61").Value2 = 1463.3334
$ws.Range("K61").Value2 = 912.4375
$ws.Range("L61").Value2 = 1463.3334
$ws.Range("M61").Value2 = -700.4375
$ws.Range("N61").Value2 = -1887.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 6968.696
$ws.Range("I74").Value2 = 975.46155
$ws.Range("K74").Value2 = 975.46155
$ws.Range("M74").Value2 = -101.46155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value2 = 6968.696
$ws.Range("I77").Value2 = 975.46155
$ws.Range("K77").Value2 = 4877.30775
$ws.Range("M77").Value2 = -509.3077499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value2 = 913.1667
$ws.Range("I122").Value2 = 790.41174
$ws.Range("J122").Value2 = 3000
$ws.Range("K122").Value2 = 2371.23522
$ws.Range("L122").Value2 = 9000
$ws.Range("M122").Value2 = 78.76477999999997
$ws.Range("N122").Value2 = -13900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value2 = 1790.2245
$ws.Range("I132").Value2 = 1681.1666
$ws.Range("J132").Value2 = 2444.5715
$ws.Range("K132").Value2 = 5043.4998
$ws.Range("L132").Value2 = 7333.7145
$ws.Range("M132").Value2 = -2513.4998
$ws.Range("N132").Value2 = -12393.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value2 = 1062.6818
$ws.Range("I136").Value2 = 912.4375
$ws.Range("J136").Value2 = 1463.3334
$ws.Range("K136").Value2 = 2737.3125
$ws.Range("L136").Value2 = 4390.0002
$ws.Range("M136").Value2 = -187.3125
$ws.Range("N136").Value2 = -9490.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 1892.9166
$ws.Range("I20").Value2 = 2130.8572
$ws.Range("J20").Value2 = 1559.8
$ws.Range("K20").Value2 = 2130.8572
$ws.Range("L20").Value2 = 1559.8
$ws.Range("M20").Value2 = -1883.8572
$ws.Range("N20").Value2 = -2053.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value2 = 20220
$ws.Range("J74").Value2 = 20220
$ws.Range("L74").Value2 = 20220
$ws.Range("N74").Value2 = -22092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H77").Value2 = 20220
$ws.Range("J77").Value2 = 20220
$ws.Range("L77").Value2 = 60660
$ws.Range("N77").Value2 = -70020

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value2 = 1480
$ws.Range("I107").Value2 = 1476.1111
$ws.Range("J107").Value2 = 1497.5
$ws.Range("K107").Value2 = 1476.1111
$ws.Range("L107").Value2 = 1497.5
$ws.Range("M107").Value2 = 443.8888999999999
$ws.Range("N107").Value2 = -5337.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 17521.268
$ws.Range("I134").Value2 = 18066.828
$ws.Range("J134").Value2 = 1700
$ws.Range("K134").Value2 = 54200.484
$ws.Range("L134").Value2 = 5100
$ws.Range("M134").Value2 = -51665.484
$ws.Range("N134").Value2 = -10170

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1728.2572
$ws.Range("I31").Value2 = 1451.9048
$ws.Range("J31").Value2 = 2142.7856
$ws.Range("K31").Value2 = 1451.9048
$ws.Range("L31").Value2 = 2142.7856
$ws.Range("M31").Value2 = -1156.9048
$ws.Range("N31").Value2 = -2732.7856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value2 = 1728.2572
$ws.Range("I34").Value2 = 1451.9048
$ws.Range("J34").Value2 = 2142.7856
$ws.Range("K34").Value2 = 1451.9048
$ws.Range("L34").Value2 = 2142.7856
$ws.Range("M34").Value2 = -1249.9048
$ws.Range("N34").Value2 = -2546.7856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value2 = 1996.7142
$ws.Range("I58").Value2 = 987
$ws.Range("J58").Value2 = 2754
$ws.Range("K58").Value2 = 987
$ws.Range("L58").Value2 = 2754
$ws.Range("M58").Value2 = -784
$ws.Range("N58").Value2 = -3160

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value2 = 1947.2858
$ws.Range("I107").Value2 = 2542.2
$ws.Range("J107").Value2 = 460
$ws.Range("K107").Value2 = 2542.2
$ws.Range("L107").Value2 = 460
$ws.Range("M107").Value2 = -622.1999999999998
$ws.Range("N107").Value2 = -4300

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value2 = 1887.5641
$ws.Range("I132").Value2 = 1760.4572
$ws.Range("K132").Value2 = 5281.3716
$ws.Range("M132").Value2 = -2751.3716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value2 = 2087.2632
$ws.Range("I134").Value2 = 2271.1516
$ws.Range("J134").Value2 = 873.6
$ws.Range("K134").Value2 = 6813.4548
$ws.Range("L134").Value2 = 2620.8
$ws.Range("M134").Value2 = -4278.4548
$ws.Range("N134").Value2 = -7690.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value2 = 1996.7142
$ws.Range("I136").Value2 = 987
$ws.Range("J136").Value2 = 2754
$ws.Range("K136").Value2 = 2961
$ws.Range("L136").Value2 = 8262
$ws.Range("M136").Value2 = -411
$ws.Range("N136").Value2 = -13362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value2 = 99.82353000000001
$ws.Range("I38").Value2 = 43.75
$ws.Range("J38").Value2 = 149.66667
$ws.Range("K38").Value2 = 131.25
$ws.Range("L38").Value2 = 449.00001
$ws.Range("M38").Value2 = 215.75
$ws.Range("N38").Value2 = -1143.00001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value2 = 675.58826
$ws.Range("J113").Value2 = 663.9286
$ws.Range("L113").Value2 = 1991.7858
$ws.Range("N113").Value2 = -6331.7858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value2 = 10000
$ws.Range("I34").Value2 = 10000
$ws.Range("K34").Value2 = 10000
$ws.Range("M34").Value2 = -9732

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value2 = 0
$ws.Range("I44").Value2 = 0
$ws.Range("J44").Value2 = 0
$ws.Range("K44").Value2 = 0
$ws.Range("L44").Value2 = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H76").Value2 = 10000
$ws.Range("I76").Value2 = 10000
$ws.Range("K76").Value2 = 10000
$ws.Range("M76").Value2 = -9685

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H79").Value2 = 10000
$ws.Range("I79").Value2 = 10000
$ws.Range("K79").Value2 = 10000
$ws.Range("M79").Value2 = -8908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value2 = 2005
$ws.Range("I126").Value2 = 1762.4
$ws.Range("K126").Value2 = 5287.200000000001
$ws.Range("M126").Value2 = -2817.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 5159.2583
$ws.Range("I132").Value2 = 5660.864
$ws.Range("K132").Value2 = 16982.592
$ws.Range("M132").Value2 = -14452.592

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value2 = 2347.2222
$ws.Range("I136").Value2 = 1415
$ws.Range("J136").Value2 = 3512.5
$ws.Range("K136").Value2 = 4245
$ws.Range("L136").Value2 = 10537.5
$ws.Range("M136").Value2 = -1695
$ws.Range("N136").Value2 = -15637.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 446.11765
$ws.Range("I107").Value2 = 198.85715
$ws.Range("J107").Value2 = 1600
$ws.Range("K107").Value2 = 596.5714499999999
$ws.Range("L107").Value2 = 4800
$ws.Range("M107").Value2 = 1323.42855
$ws.Range("N107").Value2 = -8640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 2508.0244
$ws.Range("I132").Value2 = 2333.2703
$ws.Range("J132").Value2 = 4124.5
$ws.Range("K132").Value2 = 6999.8109
$ws.Range("L132").Value2 = 12373.5
$ws.Range("M132").Value2 = -4469.8109
$ws.Range("N132").Value2 = -17433.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value2 = 915.89655
$ws.Range("I136").Value2 = 680.3
$ws.Range("J136").Value2 = 1439.4445
$ws.Range("K136").Value2 = 2040.9
$ws.Range("L136").Value2 = 4318.333500000001
$ws.Range("M136").Value2 = 509.1000000000001
$ws.Range("N136").Value2 = -9418.333500000001
